$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23; this shifts existing rows 23..138 down to 24..139
# and automatically inherits formatting (e.g. the date-style on column D) from the row below.
$ws.Rows("23").Insert()

# Populate the newly inserted row 23 with the new weekly data point.
$ws.Cells.Item(23, 1).Value = 6
$ws.Cells.Item(23, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(23, 3).Value = "Metropolitana"
$ws.Cells.Item(23, 4).Value = 44565
$ws.Cells.Item(23, 5).Value = 13
$ws.Cells.Item(23, 6).Value = 100112029
$ws.Cells.Item(23, 7).Value = "Orégano"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 35
$ws.Cells.Item(23, 11).Value = 10000
$ws.Cells.Item(23, 12).Value = 11000
$ws.Cells.Item(23, 13).Value = 10543
$ws.Cells.Item(23, 14).Value = "$/docena de atados"
$ws.Cells.Item(23, 15).Value = "Región Metropolitana"
$ws.Cells.Item(23, 16).Value = 3514
$ws.Cells.Item(23, 17).Value = 3
$ws.Cells.Item(23, 18).Value = "Hortaliza"
